{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Change: \"Expand your recruitment funnel \" -> \"Scale your recruitment funnel \"\n// i.e. the leading word \"Expand\" is replaced by \"Scale\". As a side effect of\n// this being the most-recently-edited spot, Word's \"_GoBack\" bookmark (which\n// previously sat right after \"aspiring students\" in the Events section) is\n// relocated to sit immediately after the newly typed word \"Scale\".\n\n// 1) Locate the word \"Expand\" (unique in the document) that starts the\n//    \"Expand your recruitment funnel \" run.\nlet results = context.document.body.search(\"Expand\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nlet expandRange = results.items[0];\n\n// 2) Split the run at the boundary right after \"Expand\" *before* touching its\n//    text, and do it by moving the \"_GoBack\" bookmark there. Removing the old\n//    bookmark first means Word doesn't end up with two same-named bookmarks.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst boundary = expandRange.getRange(\"End\");\nboundary.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-acquire the (now isolated, single-word) \"Expand\" range and swap it for\n//    \"Scale\": insert the new word immediately before it, then delete the old\n//    one. This keeps the edit scoped to just that run so the untouched\n//    \"to include a global audience through RALLY\" run is left exactly as-is.\nresults = context.document.body.search(\"Expand\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nexpandRange = results.items[0];\n\nexpandRange.insertText(\"Scale\", \"Before\");\nawait context.sync();\n\nresults = context.document.body.search(\"Expand\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nexpandRange = results.items[0];\nexpandRange.delete();\nawait context.sync();\n", "ps1": "# Word COM interop script\n# Change: \"Expand your recruitment funnel \" -> \"Scale your recruitment funnel \"\n# (the word \"Expand\" is replaced by \"Scale\"). As a side effect of this being\n# the most-recently-edited spot, Word's \"_GoBack\" bookmark (which previously\n# sat after \"aspiring students\" in the Events section) is relocated to sit\n# immediately after the newly typed word \"Scale\".\n\n$d = $word.ActiveDocument\n\n# 1) Locate the word \"Expand\" (unique in the document) that starts the\n#    \"Expand your recruitment funnel \" run.\n$find = $d.Content.Find\n$find.Text = \"Expand\"\n$null = $find.Execute()\n$expandRange = $find.Parent\n\n# 2) Before changing the text, mark the boundary right after \"Expand\" so the\n#    new \"_GoBack\" bookmark lands exactly between \"Scale\" and the rest of the\n#    sentence (this also splits the run there, matching how Word itself\n#    leaves the edit point as its own run boundary).\n$boundary = $expandRange.Duplicate\n$boundary.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $boundary)\n\n# 3) Now perform the actual text replacement: \"Expand\" -> \"Scale\". Re-find it\n#    since ranges can shift after the bookmark insertion.\n$find2 = $d.Content.Find\n$find2.Text = \"Expand\"\n$null = $find2.Execute()\n$find2.Parent.Text = \"Scale\"\n"}
